# MapleStory Monthly Honorable Rock Record
# Adds the "December 2019" sheet (copied from "November 2019" so that
# number formats / styles / merged cells / page setup all carry over),
# refreshes its guild-name / contribution figures, appends the new
# "HeavenSent" shared string, and re-applies the B4:C50 range selection
# that the diff shows being (re)selected across the tabs.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Duplicate "November 2019" -> new sheet placed right after it, then
#    rename it to "December 2019". Copying (rather than Worksheets.Add)
#    preserves every style index, merged cell, column width, page setup
#    and header/footer from the template sheet.
# ---------------------------------------------------------------------
$nov = $wb.Worksheets.Item("November 2019")
$nov.Copy($null, $nov)
$dec = $wb.Worksheets.Item($wb.Worksheets.Count)
$dec.Name = "December 2019"

# Month heading formula in row 2 (">    ""December 2019"":")
$dec.Range("E2").Formula = '="    """&"December 2019"&""""&":"'

# ---------------------------------------------------------------------
# 2. Overwrite the 50 guild rows (rank stays 1-50 in column B; only the
#    guild name in C and the contribution number in D change for
#    December's results - column E's formula recalculates on its own).
# ---------------------------------------------------------------------
$dec.Range("C4").Value = "Eternal"
$dec.Range("D4").Value = 240789974
$dec.Range("C5").Value = "Smile"
$dec.Range("D5").Value = 230242820
$dec.Range("C6").Value = "Savages"
$dec.Range("D6").Value = 216089993
$dec.Range("C7").Value = "Elite"
$dec.Range("D7").Value = 205319584
$dec.Range("C8").Value = "Bounce"
$dec.Range("D8").Value = 204399710
$dec.Range("C9").Value = "Spring"
$dec.Range("D9").Value = 145824555
$dec.Range("C10").Value = "Sunset"
$dec.Range("D10").Value = 136931011
$dec.Range("C11").Value = "Epic"
$dec.Range("D11").Value = 126046520
$dec.Range("C12").Value = "Beaters"
$dec.Range("D12").Value = 118695090
$dec.Range("C13").Value = "Downtime"
$dec.Range("D13").Value = 116795406
$dec.Range("C14").Value = "RainSong"
$dec.Range("D14").Value = 99533446
$dec.Range("C15").Value = "lolicafe"
$dec.Range("D15").Value = 98658827
$dec.Range("C16").Value = "Imperium"
$dec.Range("D16").Value = 97872638
$dec.Range("C17").Value = "Gintama"
$dec.Range("D17").Value = 87349880
$dec.Range("C18").Value = "Undertale"
$dec.Range("D18").Value = 86638324
$dec.Range("C19").Value = "Maha"
$dec.Range("D19").Value = 85337314
$dec.Range("C20").Value = "Remorse"
$dec.Range("D20").Value = 85244997
$dec.Range("C21").Value = "Erda"
$dec.Range("D21").Value = 81758022
$dec.Range("C22").Value = "Cleanse"
$dec.Range("D22").Value = 80272929
$dec.Range("C23").Value = "Atelier"
$dec.Range("D23").Value = 79019826
$dec.Range("C24").Value = "Broke"
$dec.Range("D24").Value = 75775892
$dec.Range("C25").Value = "Sora"
$dec.Range("D25").Value = 75631668
$dec.Range("C26").Value = "Tama"
$dec.Range("D26").Value = 74979459
$dec.Range("C27").Value = "Lithe"
$dec.Range("D27").Value = 69131473
$dec.Range("C28").Value = "Oceania"
$dec.Range("D28").Value = 69012616
$dec.Range("C29").Value = "Revive"
$dec.Range("D29").Value = 68318091
$dec.Range("C30").Value = "Rising"
$dec.Range("D30").Value = 64206429
$dec.Range("C31").Value = "Ravers"
$dec.Range("D31").Value = 62783449
$dec.Range("C32").Value = "Sugar"
$dec.Range("D32").Value = 62760688
$dec.Range("C33").Value = "Artifacts"
$dec.Range("D33").Value = 61609928
$dec.Range("C34").Value = "Fabled"
$dec.Range("D34").Value = 56077460
$dec.Range("C35").Value = "Aloe"
$dec.Range("D35").Value = 55086052
$dec.Range("C36").Value = "Earnest"
$dec.Range("D36").Value = 54168958
$dec.Range("C37").Value = "Skyfall"
$dec.Range("D37").Value = 51007455
$dec.Range("C38").Value = "CyberThreat"
$dec.Range("D38").Value = 49821525
$dec.Range("C39").Value = "Mystical"
$dec.Range("D39").Value = 49764060
$dec.Range("C40").Value = "chigga"
$dec.Range("D40").Value = 48676027
$dec.Range("C41").Value = "Fandom"
$dec.Range("D41").Value = 47515495
$dec.Range("C42").Value = "Howl"
$dec.Range("D42").Value = 46461256
$dec.Range("C43").Value = "RainDrop"
$dec.Range("D43").Value = 45174624
$dec.Range("C44").Value = "Path"
$dec.Range("D44").Value = 44934426
$dec.Range("C45").Value = "Bubbles"
$dec.Range("D45").Value = 42810639
$dec.Range("C46").Value = "Coffee"
$dec.Range("D46").Value = 42294515
$dec.Range("C47").Value = "Weibo"
$dec.Range("D47").Value = 41908988
$dec.Range("C48").Value = "Kingdom"
$dec.Range("D48").Value = 41757438
$dec.Range("C49").Value = "Exorcist"
$dec.Range("D49").Value = 41161676
$dec.Range("C50").Value = "Reboot"
$dec.Range("D50").Value = 40474527
$dec.Range("C51").Value = "Comity"
$dec.Range("D51").Value = 40296928
$dec.Range("C52").Value = "Faction"
$dec.Range("D52").Value = 39639099
$dec.Range("C53").Value = "HeavenSent"
$dec.Range("D53").Value = 37046670

# ---------------------------------------------------------------------
# 3. Re-apply the B4:C50 range selection on every sheet (picked up as
#    part of the same edit in the diff). Each sheet is activated in
#    turn and the range re-selected on it; December ends up selected
#    last so it stays the workbook's active tab.
# ---------------------------------------------------------------------
$aug = $wb.Worksheets.Item("August 2019")
$aug.Select()
$aug.Range("B4:C50").Select()

$sep = $wb.Worksheets.Item("September 2019")
$sep.Select()
$sep.Range("B4:C50").Select()

$oct_ = $wb.Worksheets.Item("October 2019")
$oct_.Select()
$oct_.Range("B4:C50").Select()

$nov.Select()
$nov.Range("B4:C50").Select()

$dec.Select()
$dec.Range("B4:C50").Select()
